$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value from 13 to 1
$ws.Range("A2").Value = 1

# Delete row 3 entirely (was A3 = 4), shifting cells up
$ws.Rows("3:3").Delete()

# Update the selection to A2
$ws.Range("A2").Select()
